$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Absence Prof")

# Fix incorrect name value in row 5 (Nom column)
$ws.Range("A5").Value = "rober "
